$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => @{ D = <new price text>; E = <new volume text> } (only keys present when that column changed)
$updates = @{
    2 = @{ D='42.088.20'; E='  -0.51%  ' }
    3 = @{ D='2.235.82'; E='  -0.45%  ' }
    4 = @{ E='  -0.01%  ' }
    5 = @{ D='249.05'; E='  +6.55%  ' }
    6 = @{ D='0.632' }
    7 = @{ D='71.90'; E='  +3.91%  ' }
    8 = @{ E='  -0.02%  ' }
    9 = @{ D='0.596'; E='  +6.04%  ' }
    10 = @{ D='41.55'; E='  +17.29%  ' }
    11 = @{ D='0.0979'; E='  -0.62%  ' }
    12 = @{ D='58.05'; E='  -0.02%  ' }
    13 = @{ D='7.19'; E='  +7.04%  ' }
    14 = @{ D='0.106'; E='  +0.01%  ' }
    15 = @{ D='2.564.02'; E='  -0.40%  ' }
    16 = @{ D='15.06'; E='  +0.79%  ' }
    17 = @{ D='0.868'; E='  +1.73%  ' }
    18 = @{ D='2.226.56'; E='  -0.52%  ' }
    19 = @{ D='41.874.68'; E='  -0.40%  ' }
    20 = @{ D='0.0₃0975'; E='  +0.01%  ' }
    21 = @{ D='6.26'; E='  +0.14%  ' }
    22 = @{ D='73.25'; E='  +0.15%  ' }
    23 = @{ D='236.65'; E='  +0.39%  ' }
    24 = @{ D='2.16'; E='  +10.39%  ' }
    25 = @{ D='4.00'; E='  +9.32%  ' }
    26 = @{ E='  +0.04%  ' }
    27 = @{ D='2.53'; E='  +7.87%  ' }
    28 = @{ D='10.76'; E='  +7.36%  ' }
    29 = @{ D='171.82'; E='  +1.86%  ' }
    30 = @{ E='  -3.23%  ' }
    31 = @{ D='20.93'; E='  +1.82%  ' }
    32 = @{ D='0.124'; E='  +3.12%  ' }
    33 = @{ D='0.126'; E='  -0.78%  ' }
    34 = @{ D='5.57'; E='  +4.76%  ' }
    35 = @{ D='0.0731'; E='  +2.03%  ' }
    36 = @{ D='4.74'; E='  +0.82%  ' }
    37 = @{ D='26.39'; E='  +24.71%  ' }
    38 = @{ D='4.01'; E='  +11.02%  ' }
    39 = @{ E='  +12.62%  ' }
    40 = @{ D='2.30'; E='  +2.11%  ' }
    41 = @{ D='6.03'; E='  +1.31%  ' }
    42 = @{ D='68.26'; E='  +3.47%  ' }
    43 = @{ D='12.17'; E='  +23.16%  ' }
    44 = @{ E='  +10.92%  ' }
    45 = @{ D='4.92'; E='  +0.42%  ' }
    46 = @{ D='8.79'; E='  -0.65%  ' }
    47 = @{ D='0.103'; E='  +1.13%  ' }
    48 = @{ D='4.71'; E='  +8.02%  ' }
    50 = @{ D='1.18'; E='  +8.43%  ' }
    51 = @{ D='1.20'; E='  +1.74%  ' }
}

foreach ($rowNum in $updates.Keys) {
    $u = $updates[$rowNum]
    if ($u.ContainsKey("D")) {
        # Force text so numeric-looking price strings (e.g. "249.05") are not
        # auto-coerced to a Number by Excel's value-assignment type inference;
        # resetting the Style afterwards drops the quote-prefix formatting so the
        # cell keeps its original (unstyled) look.
        $cell = $ws.Cells.Item($rowNum, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($rowNum, 5).Value = $u.E
    }
}
